$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain numeric "quantity" cells in column C ---
$ws.Range("C8").Value  = 29
$ws.Range("C9").Value  = 70
$ws.Range("C10").Value = 85
$ws.Range("C11").Value = 5
$ws.Range("C12").Value = 34
$ws.Range("C13").Value = 52
$ws.Range("C14").Value = 18
$ws.Range("C15").Value = 25
$ws.Range("C16").Value = 34
$ws.Range("C17").Value = 68

# --- Amount cells that must stay text (e.g. "17920.00") ---
# Plain .Value assignment of a numeric-looking string gets auto-converted
# to a number by Excel, dropping the literal formatting ("17920.00" -> 17920).
# Route the new text through a formula + paste-values-only round trip so the
# cells keep their original text ("str") cell type without picking up any
# new number-format/style (which a leading-apostrophe or NumberFormat="@"
# trick would otherwise introduce).
$ws.Range("G9").Formula  = '="17920.00"'
$ws.Range("G10").Formula = '="40120.00"'
$ws.Range("G11").Formula = '="3310.00"'
$ws.Range("G13").Formula = '="7072.00"'
$ws.Range("G14").Formula = '="414.00"'
$ws.Range("G19").Formula = '="68836.00"'
$ws.Range("H19").Formula = '="68836.00"'
$ws.Range("G21").Formula = '="68836.00"'
$ws.Range("H21").Formula = '="68836.00"'

$amountRanges = @("G9:G11", "G13:G14", "G19:H19", "G21:H21")
foreach ($addr in $amountRanges) {
    $rng = $ws.Range($addr)
    $rng.Copy()
    $rng.PasteSpecial(-4163)
}
